# Update the "想去人数" (want-to-go count) figures in column F, rows 2-6,
# on both the "展览" and "全部类型" worksheets, reflecting freshly
# regenerated data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 6538
    3 = 38
    4 = 191
    5 = 1022
    6 = 124
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
